$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rename / re-word a handful of existing (shared) category labels.
#    Each of these strings is reused across many rows in the sheet, so
#    a sheet-wide exact-match replace mirrors what happened upstream.
#    xlWhole = 1, xlByRows = 1 (SearchOrder)
# ------------------------------------------------------------------
$xlWhole = 1
$xlByRows = 1
$ws.Cells.Replace("NSFW Discord Server Invite", "NSFW Discord Invite", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("Malicious discord server", "Discord Servers", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("Discord accounts", "Discord Accounts", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("Cloned Steam login page", "Cloned Steam Pages", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("Funpay scam", "Funpay Fraud", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("Funpay offer", "Funpay Offers", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("Fake Discord Nitro gift", "Fake Nitro Gift", $xlWhole, $xlByRows, $false, $false, $false)
$ws.Cells.Replace("Browser Session Hijacking", "Session Hijacking", $xlWhole, $xlByRows, $false, $false, $false)

# ------------------------------------------------------------------
# 2) Append two new case rows (664, 665), cloned from the formatting
#    of the previous row (663) which is the same "Steam phishing via
#    Markdown URL Obfuscation" case pattern.
# ------------------------------------------------------------------
$ws.Range("A663:L663").Copy()
$ws.Range("A664:L665").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rows = @(664, 665)
$discordIds = @("1138023535651074070", "691835106565947473")
$usernames  = @("e50ausf.m", "blizarice_")
$links      = @("https://steanmnscommunity.com/105381409", "https://u.to/h-UkIg")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 2).Value = "2025-03-12"
    $ws.Cells.Item($r, 3).Value = $discordIds[$i]
    $ws.Cells.Item($r, 4).Value = $usernames[$i]
    $ws.Cells.Item($r, 5).Value = "Markdown URL Obfuscation"
    $ws.Cells.Item($r, 6).Value = "Phishing site"
    $ws.Cells.Item($r, 7).Value = "Cloned Steam Pages"
    $ws.Cells.Item($r, 8).Value = "Steam accounts"
    $ws.Cells.Item($r, 9).Value = "Steam"
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 10), $links[$i])
    $ws.Cells.Item($r, 10).Font.Underline = 2
    $ws.Cells.Item($r, 10).Font.Color = 16711680
    $ws.Cells.Item($r, 11).Value = "UNKNOWN"
    $ws.Cells.Item($r, 12).Value = "UNKNOWN"
}
